$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New technology columns: solar_th1 (D) and pvt1 (E), mirroring the
# existing net1 / CHP1 header columns (B / C).
$ws.Range("D1").Value = "solar_th1"
$ws.Range("E1").Value = "pvt1"

# Copy the header formatting (bold, bordered, centered) from the existing
# header cell C1 onto the two new header cells so they share the same style.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new data rows with zeros, matching B2:C3 / B3:C3.
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
